# [Minno, Juan] #203 fixed tests for import students
# The header column previously labelled "Name" was renamed to the more
# descriptive "Name of the Student" on both test sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("happy-path")
$ws1.Range("B4").Value = "Name of the Student"

$ws2 = $wb.Worksheets.Item("offset-columns")
$ws2.Range("C4").Value = "Name of the Student"

# Leave the cursor where the edits were made, matching the author's
# last-saved selection state on each sheet.
[void]$ws1.Range("I4").Select()

$ws2.Activate()
[void]$ws2.Range("C4").Select()
